# Consolidate some of the points for delineating catchments
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1) Remove the extraneous "Sheet1" worksheet (a filtered helper/staging sheet).
#    This also removes its _xlnm._FilterDatabase defined name automatically.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

# 2) On "Summary_all_points", merge the "Tanaelva ved Jalvvivárri" (Vannmiljø / Water chem)
#    row with the "Tana v/Storfossen" (NVE / SS) row into a single combined row.
$ws = $wb.Worksheets.Item("Summary_all_points")

# Row 28 ("234-90042" / "Tanaelva ved Jalvvivárri" / Vannmiljø / Water chem) is dropped;
# deleting it shifts row 29 ("Polmak") up to row 28, and row 30 ("Tana v/Storfossen") up to row 29.
$ws.Rows.Item(28).Delete()

# Now row 29 holds the old "Tana v/Storfossen" / NVE / SS data - update it to reflect the
# consolidated station that combines both former rows.
$ws.Range("E29").Value = "NVE, Vannmiljø"
$ws.Range("F29").Value = "SS, water chem"
$ws.Range("B29").Value = "Tana v/Storfossen & Jalvvivárri "

# 3) Refresh the AutoFilter range now that the sheet has one fewer row (36 -> 35).
$ws.AutoFilterMode = $false
$ws.Range("A1:F35").AutoFilter()

# 4) Keep the workbook-level _FilterDatabase defined name for this sheet in sync too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Summary_all_points!_FilterDatabase") {
        $n.RefersTo = "=Summary_all_points!`$A`$1:`$F`$35"
    }
}

$wb.Save()
